$wb = $excel.ActiveWorkbook

# Row 6 on ALC (diff hunk -932,22 +932,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value2 = 75
$ws.Cells.Item(6, 9).Value2 = 75
$ws.Cells.Item(6, 11).Value2 = 225
$ws.Cells.Item(6, 13).Value2 = -113

# Row 8 on ALC (diff hunk -1030,19 +1030,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value2 = 100
$ws.Cells.Item(8, 9).Value2 = 100
$ws.Cells.Item(8, 11).Value2 = 300
$ws.Cells.Item(8, 13).Value2 = -161

# Row 31 on ALC (diff hunk -2157,22 +2160,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(31, 8).Value2 = 2346.3333
$ws.Cells.Item(31, 9).Value2 = 2346.3333
$ws.Cells.Item(31, 11).Value2 = 7038.999899999999
$ws.Cells.Item(31, 13).Value2 = -6808.999899999999

# Row 116 on ALC (diff hunk -6400,22 +6403,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 9).Value2 = 11114580
$ws.Cells.Item(116, 10).Value2 = 4480.8
$ws.Cells.Item(116, 11).Value2 = 11114580
$ws.Cells.Item(116, 12).Value2 = 4480.8
$ws.Cells.Item(116, 13).Value2 = -11111138
$ws.Cells.Item(116, 14).Value2 = -11364.8

# Row 129 on ALC (diff hunk -7037,22 +7040,22)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129, 8).Value2 = 1012.125
$ws.Cells.Item(129, 9).Value2 = 871.1429000000001
$ws.Cells.Item(129, 11).Value2 = 2613.4287
$ws.Cells.Item(129, 13).Value2 = 2386.5713

# Row 139 on ALC (diff hunk -7536,25 +7539,25)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(139, 8).Value2 = 115756
$ws.Cells.Item(139, 10).Value2 = 100780
$ws.Cells.Item(139, 12).Value2 = 100780
$ws.Cells.Item(139, 14).Value2 = -111060

# Row 140 on ALC (diff hunk -7588,25 +7591,25)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(140, 8).Value2 = 172166.17
$ws.Cells.Item(140, 10).Value2 = 218249.25
$ws.Cells.Item(140, 12).Value2 = 218249.25
$ws.Cells.Item(140, 14).Value2 = -228609.25

# Row 11 on ARM (diff hunk -8242,22 +8245,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(11, 8).Value2 = 25002000
$ws.Cells.Item(11, 9).Value2 = 25002000
$ws.Cells.Item(11, 11).Value2 = 25002000
$ws.Cells.Item(11, 13).Value2 = -25001856

# Row 32 on ARM (diff hunk -9280,25 +9283,25)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 1885583.1
$ws.Cells.Item(32, 9).Value2 = 1932.6545
$ws.Cells.Item(32, 10).Value2 = 27785778
$ws.Cells.Item(32, 11).Value2 = 1932.6545
$ws.Cells.Item(32, 12).Value2 = 27785778
$ws.Cells.Item(32, 13).Value2 = -1645.6545
$ws.Cells.Item(32, 14).Value2 = -27786352

# Row 61 on ARM (diff hunk -10704,22 +10707,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value2 = 5529.6304
$ws.Cells.Item(61, 9).Value2 = 12200
$ws.Cells.Item(61, 11).Value2 = 12200
$ws.Cells.Item(61, 13).Value2 = -11988

# Row 74 on ARM (diff hunk -11320,22 +11323,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value2 = 4813.684
$ws.Cells.Item(74, 9).Value2 = 5930.636
$ws.Cells.Item(74, 11).Value2 = 5930.636
$ws.Cells.Item(74, 13).Value2 = -5056.636

# Row 77 on ARM (diff hunk -11467,22 +11470,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value2 = 4813.684
$ws.Cells.Item(77, 9).Value2 = 5930.636
$ws.Cells.Item(77, 11).Value2 = 29653.18
$ws.Cells.Item(77, 13).Value2 = -25285.18

# Row 132 on ARM (diff hunk -14099,22 +14102,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value2 = 913560.0600000001
$ws.Cells.Item(132, 9).Value2 = 1115538.4
$ws.Cells.Item(132, 11).Value2 = 3346615.2
$ws.Cells.Item(132, 13).Value2 = -3344085.2

# Row 133 on ARM (diff hunk -14151,22 +14154,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(133, 8).Value2 = 219632.33
$ws.Cells.Item(133, 10).Value2 = 219632.33
$ws.Cells.Item(133, 12).Value2 = 219632.33
$ws.Cells.Item(133, 14).Value2 = -224692.33

# Row 136 on ARM (diff hunk -14292,22 +14295,22)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value2 = 5529.6304
$ws.Cells.Item(136, 9).Value2 = 12200
$ws.Cells.Item(136, 11).Value2 = 36600
$ws.Cells.Item(136, 13).Value2 = -34050

# Row 141 on ARM (diff hunk -14534,22 +14537,19)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(141, 8).Value2 = 39000
$ws.Cells.Item(141, 9).Value2 = 0
$ws.Cells.Item(141, 11).Value2 = 0
$ws.Cells.Item(141, 13).ClearContents()

# Row 80 on BSM (diff hunk -18466,22 +18466,22)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value2 = 13344463
$ws.Cells.Item(80, 9).Value2 = 1852.8334
$ws.Cells.Item(80, 11).Value2 = 1852.8334
$ws.Cells.Item(80, 13).Value2 = -854.8334

# Row 83 on BSM (diff hunk -18616,22 +18616,22)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(83, 8).Value2 = 13344463
$ws.Cells.Item(83, 9).Value2 = 1852.8334
$ws.Cells.Item(83, 11).Value2 = 9264.166999999999
$ws.Cells.Item(83, 13).Value2 = -4272.166999999999

# Row 94 on BSM (diff hunk -19146,22 +19146,22)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value2 = 36027.71
$ws.Cells.Item(94, 9).Value2 = 1417.3684
$ws.Cells.Item(94, 11).Value2 = 1417.3684
$ws.Cells.Item(94, 13).Value2 = -966.3684000000001

# Row 134 on BSM (diff hunk -21061,22 +21061,22)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value2 = 1321809.4
$ws.Cells.Item(134, 9).Value2 = 1618113.1
$ws.Cells.Item(134, 11).Value2 = 4854339.300000001
$ws.Cells.Item(134, 13).Value2 = -4851804.300000001

# Row 31 on CRP (diff hunk -22995,25 +22995,25)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value2 = 5626.625
$ws.Cells.Item(31, 9).Value2 = 1321.2858
$ws.Cells.Item(31, 10).Value2 = 8975.223
$ws.Cells.Item(31, 11).Value2 = 1321.2858
$ws.Cells.Item(31, 12).Value2 = 8975.223
$ws.Cells.Item(31, 13).Value2 = -1026.2858
$ws.Cells.Item(31, 14).Value2 = -9565.223

# Row 34 on CRP (diff hunk -23148,25 +23148,25)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value2 = 5626.625
$ws.Cells.Item(34, 9).Value2 = 1321.2858
$ws.Cells.Item(34, 10).Value2 = 8975.223
$ws.Cells.Item(34, 11).Value2 = 1321.2858
$ws.Cells.Item(34, 12).Value2 = 8975.223
$ws.Cells.Item(34, 13).Value2 = -1119.2858
$ws.Cells.Item(34, 14).Value2 = -9379.223

# Row 122 on CRP (diff hunk -27406,22 +27406,22)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value2 = 6385.5
$ws.Cells.Item(122, 9).Value2 = 3369.5
$ws.Cells.Item(122, 11).Value2 = 10108.5
$ws.Cells.Item(122, 13).Value2 = -7658.5

# Row 6 on CUL (diff hunk -28682,22 +28682,19)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value2 = 0
$ws.Cells.Item(6, 9).Value2 = 0
$ws.Cells.Item(6, 11).Value2 = 0
$ws.Cells.Item(6, 13).ClearContents()

# Row 101 on CUL (diff hunk -33475,25 +33472,25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(101, 8).Value2 = 27002.2
$ws.Cells.Item(101, 10).Value2 = 28332.889
$ws.Cells.Item(101, 12).Value2 = 84998.667
$ws.Cells.Item(101, 14).Value2 = -89866.667

# Row 124 on CUL (diff hunk -34635,25 +34632,25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(124, 8).Value2 = 5910.2
$ws.Cells.Item(124, 9).Value2 = 4749
$ws.Cells.Item(124, 10).Value2 = 10555
$ws.Cells.Item(124, 11).Value2 = 14247
$ws.Cells.Item(124, 12).Value2 = 31665
$ws.Cells.Item(124, 13).Value2 = -9337
$ws.Cells.Item(124, 14).Value2 = -41485

# Row 131 on CUL (diff hunk -34990,25 +34987,25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value2 = 72226440
$ws.Cells.Item(131, 9).Value2 = 76195920
$ws.Cells.Item(131, 10).Value2 = 66669176
$ws.Cells.Item(131, 11).Value2 = 228587760
$ws.Cells.Item(131, 12).Value2 = 200007528
$ws.Cells.Item(131, 13).Value2 = -228582720
$ws.Cells.Item(131, 14).Value2 = -200017608

# Row 134 on CUL (diff hunk -35140,22 +35137,22)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value2 = 142883140
$ws.Cells.Item(134, 9).Value2 = 142883140
$ws.Cells.Item(134, 11).Value2 = 428649420
$ws.Cells.Item(134, 13).Value2 = -428644350

# Row 137 on CUL (diff hunk -35293,25 +35290,25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value2 = 3131.35
$ws.Cells.Item(137, 9).Value2 = 1580.6666
$ws.Cells.Item(137, 10).Value2 = 3795.9285
$ws.Cells.Item(137, 11).Value2 = 4741.9998
$ws.Cells.Item(137, 12).Value2 = 11387.7855
$ws.Cells.Item(137, 13).Value2 = 358.0002000000004
$ws.Cells.Item(137, 14).Value2 = -21587.7855

# Row 126 on GSM (diff hunk -41615,25 +41612,25)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value2 = 21748964
$ws.Cells.Item(126, 10).Value2 = 17333.111
$ws.Cells.Item(126, 12).Value2 = 51999.333
$ws.Cells.Item(126, 14).Value2 = -56939.333

# Row 46 on LTW (diff hunk -44685,25 +44682,25)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value2 = 38463890
$ws.Cells.Item(46, 9).Value2 = 1378.75
$ws.Cells.Item(46, 10).Value2 = 55558340
$ws.Cells.Item(46, 11).Value2 = 1378.75
$ws.Cells.Item(46, 12).Value2 = 55558340
$ws.Cells.Item(46, 13).Value2 = -1190.75
$ws.Cells.Item(46, 14).Value2 = -55558716

# Row 55 on LTW (diff hunk -45138,25 +45135,25)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value2 = 5127.4443
$ws.Cells.Item(55, 10).Value2 = 6500.8184
$ws.Cells.Item(55, 12).Value2 = 6500.8184
$ws.Cells.Item(55, 14).Value2 = -6846.8184

# Row 93 on LTW (diff hunk -46982,25 +46979,25)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value2 = 1155.1578
$ws.Cells.Item(93, 9).Value2 = 1112.1333
$ws.Cells.Item(93, 10).Value2 = 1316.5
$ws.Cells.Item(93, 11).Value2 = 1112.1333
$ws.Cells.Item(93, 12).Value2 = 1316.5
$ws.Cells.Item(93, 13).Value2 = 135.8667
$ws.Cells.Item(93, 14).Value2 = -3812.5

# Row 94 on LTW (diff hunk -47034,22 +47031,22)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(94, 8).Value2 = 28500
$ws.Cells.Item(94, 10).Value2 = 28500
$ws.Cells.Item(94, 12).Value2 = 28500
$ws.Cells.Item(94, 14).Value2 = -29852

# Row 24 on WVR (diff hunk -50528,22 +50525,19)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(24, 8).Value2 = 0
$ws.Cells.Item(24, 10).Value2 = 0
$ws.Cells.Item(24, 12).Value2 = 0
$ws.Cells.Item(24, 14).ClearContents()

# Row 122 on WVR (diff hunk -55270,22 +55264,22)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value2 = 11283.414
$ws.Cells.Item(122, 9).Value2 = 3617.4211
$ws.Cells.Item(122, 11).Value2 = 10852.2633
$ws.Cells.Item(122, 13).Value2 = -8402.263300000001

# Row 136 on WVR (diff hunk -55956,22 +55950,22)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value2 = 62541308
$ws.Cells.Item(136, 9).Value2 = 125053550
$ws.Cells.Item(136, 11).Value2 = 375160650
$ws.Cells.Item(136, 13).Value2 = -375158100
